$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origFmt = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = $origFmt
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "69.863.53"
$ws.Range("D3").Value = "3.530.64"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue "D5" "605.68"
$ws.Range("E5").Value = "  -0.35%  "
Set-TextValue "D6" "195.94"
$ws.Range("E6").Value = "  +2.12%  "
Set-TextValue "D7" "0.625"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -4.68%  "
$ws.Range("E10").Value = "  -2.66%  "
Set-TextValue "D11" "53.46"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  -1.40%  "
Set-TextValue "D13" "9.47"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").Value = "4.092.19"
$ws.Range("E14").Value = "  +0.67%  "
Set-TextValue "D15" "597.23"
$ws.Range("E15").Value = "  -3.57%  "
$ws.Range("D16").Value = "69.977.21"
$ws.Range("E16").Value = "  +0.05%  "
Set-TextValue "D17" "12.73"
$ws.Range("E17").Value = "  +0.09%  "
Set-TextValue "D18" "19.00"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "3.527.32"
$ws.Range("E19").Value = "  +0.36%  "
Set-TextValue "D20" "0.122"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("E21").Value = "  -0.63%  "
Set-TextValue "D22" "17.92"
$ws.Range("E22").Value = "  +1.11%  "
Set-TextValue "D23" "103.23"
$ws.Range("E23").Value = "  -2.52%  "
$ws.Range("E24").Value = "  +3.12%  "
Set-TextValue "D25" "4.62"
$ws.Range("E25").Value = "  -0.59%  "
Set-TextValue "D26" "3.05"
$ws.Range("E26").Value = "  +0.77%  "
Set-TextValue "D27" "10.79"
$ws.Range("E27").Value = "  -1.91%  "
Set-TextValue "D28" "9.53"
$ws.Range("E28").Value = "  -3.20%  "
$ws.Range("E29").Value = "  -2.45%  "
Set-TextValue "D30" "7.08"
$ws.Range("E30").Value = "  -0.64%  "
Set-TextValue "D31" "4.24"
$ws.Range("E31").Value = "  +1.98%  "
Set-TextValue "D32" "12.34"
$ws.Range("E32").Value = "  -2.43%  "
$ws.Range("E33").Value = "  -0.02%  "
Set-TextValue "D34" "63.46"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").Value = "3.776.06"
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("D37").Value = "0.0₃0812"
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("E38").Value = "  +0.10%  "
Set-TextValue "D39" "508.78"
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("E40").Value = "  -0.03%  "
Set-TextValue "D41" "3.57"
$ws.Range("E41").Value = "  +0.01%  "
Set-TextValue "D42" "36.44"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("E43").Value = "  -3.15%  "
$ws.Range("E44").Value = "  -2.87%  "
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("E47").Value = "  -2.36%  "
Set-TextValue "D48" "1.00"
$ws.Range("E48").Value = "  +0.10%  "
Set-TextValue "D49" "8.47"
$ws.Range("E49").Value = "  -3.17%  "
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue "D50" "0.000249"
$ws.Range("E50").Value = "  +5.66%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D51" "1.34"
$ws.Range("E51").Value = "  +3.81%  "
